$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pollutant_concentration")

# Insert a single new row at 12 (within the "Черкаська ТЕЦ" block) for
# "Діоксид сірки" / 0.5. This pushes the old rows 12-19 down to 13-20.
$ws.Rows.Item(12).Insert()

# --- Fill in the brand-new row 12 ---
$ws.Range("A12").Value2 = "ПрАТ ""Черкаське хімволокно""  ВП ""Черкаська ТЕЦ"""
$ws.Range("B12").Value2 = "Діоксид сірки"
$ws.Range("C12").Value2 = 0.5
$ws.Rows.Item(12).RowHeight = 30

# Update selection to match the saved view state.
$ws.Range("C12").Select()
